# Auto-generated script applying market-price data refresh to H:N columns
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4719.0264
$ws.Range("I100").Value = 2947.25
$ws.Range("J100").Value = 5191.5
$ws.Range("K100").Value = 2947.25
$ws.Range("L100").Value = 5191.5
$ws.Range("M100").Value = -2406.25
$ws.Range("N100").Value = -6273.5
$ws.Range("H137").Value = 2093.8853
$ws.Range("I137").Value = 1960.6875
$ws.Range("K137").Value = 5882.0625
$ws.Range("M137").Value = -3332.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 377178.47
$ws.Range("I61").Value = 245907.83
$ws.Range("K61").Value = 245907.83
$ws.Range("M61").Value = -245695.83
$ws.Range("H74").Value = 205460.9
$ws.Range("I74").Value = 239172.1
$ws.Range("K74").Value = 239172.1
$ws.Range("M74").Value = -238298.1
$ws.Range("H77").Value = 205460.9
$ws.Range("I77").Value = 239172.1
$ws.Range("K77").Value = 1195860.5
$ws.Range("M77").Value = -1191492.5
$ws.Range("H132").Value = 3065.691
$ws.Range("I132").Value = 2913.8975
$ws.Range("J132").Value = 3435.6875
$ws.Range("K132").Value = 8741.692500000001
$ws.Range("L132").Value = 10307.0625
$ws.Range("M132").Value = -6211.692500000001
$ws.Range("N132").Value = -15367.0625
$ws.Range("H136").Value = 377178.47
$ws.Range("I136").Value = 245907.83
$ws.Range("K136").Value = 737723.49
$ws.Range("M136").Value = -735173.49

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4751.6313
$ws.Range("I86").Value = 5726.1665
$ws.Range("K86").Value = 5726.1665
$ws.Range("M86").Value = -4603.1665
$ws.Range("H89").Value = 4751.6313
$ws.Range("I89").Value = 5726.1665
$ws.Range("K89").Value = 28630.8325
$ws.Range("M89").Value = -23014.8325
$ws.Range("H97").Value = 10868.875
$ws.Range("I97").Value = 10278.714
$ws.Range("J97").Value = 15000
$ws.Range("K97").Value = 10278.714
$ws.Range("L97").Value = 15000
$ws.Range("M97").Value = -9287.714
$ws.Range("N97").Value = -16982
$ws.Range("H134").Value = 2145.3289
$ws.Range("I134").Value = 1710.322
$ws.Range("J134").Value = 3978.5715
$ws.Range("K134").Value = 5130.965999999999
$ws.Range("L134").Value = 11935.7145
$ws.Range("M134").Value = -2595.965999999999
$ws.Range("N134").Value = -17005.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2581
$ws.Range("I31").Value = 1818.449
$ws.Range("J31").Value = 4205.5654
$ws.Range("K31").Value = 1818.449
$ws.Range("L31").Value = 4205.5654
$ws.Range("M31").Value = -1523.449
$ws.Range("N31").Value = -4795.5654
$ws.Range("H34").Value = 2581
$ws.Range("I34").Value = 1818.449
$ws.Range("J34").Value = 4205.5654
$ws.Range("K34").Value = 1818.449
$ws.Range("L34").Value = 4205.5654
$ws.Range("M34").Value = -1616.449
$ws.Range("N34").Value = -4609.5654
$ws.Range("H132").Value = 2721.0435
$ws.Range("I132").Value = 1256
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 3768
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1238
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 1940.5143
$ws.Range("I134").Value = 1255.6364
$ws.Range("J134").Value = 3099.5386
$ws.Range("K134").Value = 3766.9092
$ws.Range("L134").Value = 9298.6158
$ws.Range("M134").Value = -1231.9092
$ws.Range("N134").Value = -14368.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 89
$ws.Range("I14").Value = 89
$ws.Range("K14").Value = 267
$ws.Range("M14").Value = -94
$ws.Range("H62").Value = 3165
$ws.Range("I62").Value = 2490
$ws.Range("J62").Value = 3300
$ws.Range("K62").Value = 7470
$ws.Range("L62").Value = 9900
$ws.Range("M62").Value = -6784
$ws.Range("N62").Value = -11272
$ws.Range("H65").Value = 3165
$ws.Range("I65").Value = 2490
$ws.Range("J65").Value = 3300
$ws.Range("K65").Value = 22410
$ws.Range("L65").Value = 29700
$ws.Range("M65").Value = -18978
$ws.Range("N65").Value = -36564

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 3500
$ws.Range("J33").Value = 3500
$ws.Range("L33").Value = 3500
$ws.Range("N33").Value = -4004
$ws.Range("H97").Value = 1081.3125
$ws.Range("I97").Value = 811.1111
$ws.Range("J97").Value = 1428.7142
$ws.Range("K97").Value = 811.1111
$ws.Range("L97").Value = 1428.7142
$ws.Range("M97").Value = -315.1111
$ws.Range("N97").Value = -2420.7142
$ws.Range("H122").Value = 1934.875
$ws.Range("I122").Value = 1355.8
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 4067.4
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -1617.4
$ws.Range("N122").Value = -13600
$ws.Range("H132").Value = 3345.05
$ws.Range("I132").Value = 3269.0344
$ws.Range("J132").Value = 3545.4546
$ws.Range("K132").Value = 9807.1032
$ws.Range("L132").Value = 10636.3638
$ws.Range("M132").Value = -7277.1032
$ws.Range("N132").Value = -15696.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 320.0909
$ws.Range("I55").Value = 177.625
$ws.Range("J55").Value = 700
$ws.Range("K55").Value = 177.625
$ws.Range("L55").Value = 700
$ws.Range("M55").Value = -4.625
$ws.Range("N55").Value = -1046
$ws.Range("H100").Value = 50004490
$ws.Range("I100").Value = 7473.222
$ws.Range("J100").Value = 90911140
$ws.Range("K100").Value = 7473.222
$ws.Range("L100").Value = 90911140
$ws.Range("M100").Value = -6932.222
$ws.Range("N100").Value = -90912222
$ws.Range("H132").Value = 8848.6
$ws.Range("I132").Value = 3044
$ws.Range("J132").Value = 18671.77
$ws.Range("K132").Value = 9132
$ws.Range("L132").Value = 56015.31
$ws.Range("M132").Value = -6602
$ws.Range("N132").Value = -61075.31

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2567.2856
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 3114.2
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 9342.599999999999
$ws.Range("M126").Value = -1130
$ws.Range("N126").Value = -14282.6
$ws.Range("H132").Value = 2006.075
$ws.Range("I132").Value = 1427.8077
$ws.Range("J132").Value = 3080
$ws.Range("K132").Value = 4283.4231
$ws.Range("L132").Value = 9240
$ws.Range("M132").Value = -1753.4231
$ws.Range("N132").Value = -14300
$ws.Range("H136").Value = 15320017
$ws.Range("I136").Value = 22245954
$ws.Range("J136").Value = 478724.1
$ws.Range("K136").Value = 66737862
$ws.Range("L136").Value = 1436172.3
$ws.Range("M136").Value = -66735312
$ws.Range("N136").Value = -1441272.3
